# Update the cryptocurrency price/volume snapshot (GitHub Actions refresh).
# Columns D (Price) and E (Volume(1h)) hold text-formatted values (e.g. "29.207.12",
# "  -0.95%  "), so we force Text number format before writing so Excel does not
# coerce the strings into numbers/dates, then restore the default "Normal" style
# afterwards so the cells keep no explicit style (matching the original workbook).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fmtRange = $ws.Range("D2:E51")
$fmtRange.NumberFormat = "@"

$ws.Range('D2').Value = '29.207.12'
$ws.Range('E2').Value = '  -0.95%  '
$ws.Range('D3').Value = '1.867.57'
$ws.Range('E3').Value = '  -0.46%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = '0.7091'
$ws.Range('E5').Value = '  -0.95%  '
$ws.Range('D6').Value = '241.91'
$ws.Range('E6').Value = '  +0.15%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').Value = '0.3111'
$ws.Range('E8').Value = '  +0.11%  '
$ws.Range('E9').Value = '  -3.87%  '
$ws.Range('D10').Value = '24.67'
$ws.Range('E10').Value = '  -2.68%  '
$ws.Range('D11').Value = '0.08370'
$ws.Range('E11').Value = '  +1.01%  '
$ws.Range('D12').Value = '1.863.20'
$ws.Range('E12').Value = '  -0.68%  '
$ws.Range('D13').Value = '5.228'
$ws.Range('E13').Value = '  -1.04%  '
$ws.Range('E14').Value = '  -2.70%  '
$ws.Range('D15').Value = '91.22'
$ws.Range('E15').Value = '  +0.01%  '
$ws.Range('D16').Value = '29.213.17'
$ws.Range('E16').Value = '  -0.94%  '
$ws.Range('D17').Value = '5.948'
$ws.Range('E17').Value = '  +0.25%  '
$ws.Range('D18').Value = '243.48'
$ws.Range('E18').Value = '  -0.81%  '
$ws.Range('D19').Value = '0.000007832'
$ws.Range('E19').Value = '  -0.70%  '
$ws.Range('D20').Value = '2.114.56'
$ws.Range('E20').Value = '  -0.54%  '
$ws.Range('E21').Value = '  -2.02%  '
$ws.Range('D22').Value = '0.9997'
$ws.Range('E22').Value = '  +0.01%  '
$ws.Range('D23').Value = '7.853'
$ws.Range('E23').Value = '  -1.73%  '
$ws.Range('E24').Value = '  +0.11%  '
$ws.Range('D25').Value = '0.1625'
$ws.Range('E25').Value = '  +1.00%  '
$ws.Range('E26').Value = '  -0.24%  '
$ws.Range('D27').Value = '8.959'
$ws.Range('E27').Value = '  -1.05%  '
$ws.Range('D28').Value = '18.52'
$ws.Range('E28').Value = '  +1.03%  '
$ws.Range('D29').Value = '1.509'
$ws.Range('E29').Value = '  +1.07%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').Value = '1.318'
$ws.Range('E30').Value = '  -3.02%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').Value = '4.397'
$ws.Range('E31').Value = '  +0.16%  '
$ws.Range('D32').Value = '4.242'
$ws.Range('E32').Value = '  +3.06%  '
$ws.Range('D33').Value = '0.05135'
$ws.Range('E33').Value = '  -2.51%  '
$ws.Range('D34').Value = '0.7954'
$ws.Range('E34').Value = '  +9.34%  '
$ws.Range('D35').Value = '1.910'
$ws.Range('E35').Value = '  -2.37%  '
$ws.Range('E36').Value = '  -2.75%  '
$ws.Range('E37').Value = '  +0.44%  '
$ws.Range('D38').Value = '0.01855'
$ws.Range('E38').Value = '  -0.80%  '
$ws.Range('D39').Value = '2.709'
$ws.Range('E39').Value = '  +0.03%  '
$ws.Range('D40').Value = '1.153.97'
$ws.Range('E40').Value = '  -5.57%  '
$ws.Range('D41').Value = '6.361'
$ws.Range('E41').Value = '  +3.84%  '
$ws.Range('D42').Value = '0.8961'
$ws.Range('E42').Value = '  -1.56%  '
$ws.Range('D43').Value = '73.18'
$ws.Range('E43').Value = '  -1.02%  '
$ws.Range('D44').Value = '0.9999'
$ws.Range('E44').Value = '  +0.00%  '
$ws.Range('D45').Value = '103.20'
$ws.Range('E45').Value = '  +0.84%  '
$ws.Range('D46').Value = '2.011.94'
$ws.Range('E46').Value = '  -0.11%  '
$ws.Range('D47').Value = '0.5164'
$ws.Range('E47').Value = '  -2.25%  '
$ws.Range('E48').Value = '  -1.17%  '
$ws.Range('D49').Value = '9.336'
$ws.Range('E49').Value = '  -0.05%  '
$ws.Range('E50').Value = '  -0.79%  '
$ws.Range('D51').Value = '0.4294'
$ws.Range('E51').Value = '  -0.70%  '

# Restore default styling so the text-formatted cells have no explicit style index.
$fmtRange.Style = "Normal"
